# Updated symbol list on Sun Jan  8 06:58:48 UTC 2023 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures for the crypto rows on the active worksheet. Values are written
# as text (matching the workbook's original inline-string cells) by
# temporarily forcing a text number format, then restoring the cell style
# to "Normal" so no extra formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}


Set-TextValue "D2" "261.07"
Set-TextValue "E2" "-0.03%"
Set-TextValue "D3" "26.90"
Set-TextValue "E3" "-1.55%"
Set-TextValue "D4" "4.701"
Set-TextValue "E4" "-0.24%"
Set-TextValue "D5" "0.06217"
Set-TextValue "E5" "2.31%"
Set-TextValue "D6" "6.753"
Set-TextValue "E6" "1.34%"
Set-TextValue "D7" "0.8516"
Set-TextValue "E7" "0.52%"
Set-TextValue "D8" "0.9138"
Set-TextValue "E8" "-0.23%"
Set-TextValue "D9" "0.1405"
Set-TextValue "E9" "-0.13%"
Set-TextValue "D10" "0.04937"
Set-TextValue "E10" "2.04%"
Set-TextValue "D11" "0.07068"
Set-TextValue "E11" "-0.45%"
Set-TextValue "E12" "-1.54%"
Set-TextValue "D13" "0.09047"
Set-TextValue "E13" "-0.32%"
Set-TextValue "D14" "0.001527"
Set-TextValue "E14" "-1.37%"
Set-TextValue "D15" "0.0006167"
Set-TextValue "E15" "1.73%"
Set-TextValue "D16" "0.006050"
Set-TextValue "E16" "-1.53%"
Set-TextValue "D17" "3.443"
Set-TextValue "E17" "-0.16%"
Set-TextValue "D18" "3.179"
Set-TextValue "E18" "1.00%"
Set-TextValue "D19" "2.146"
Set-TextValue "E19" "-1.39%"
Set-TextValue "D21" "0.1310"
Set-TextValue "E21" "0.78%"
Set-TextValue "D22" "4.106"
Set-TextValue "E22" "0.20%"
Set-TextValue "D24" "0.001207"
Set-TextValue "E24" "-0.94%"
Set-TextValue "D25" "0.004075"
Set-TextValue "E25" "4.17%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "0.02%"
Set-TextValue "E27" "4.39%"
Set-TextValue "D40" "0.03948"
Set-TextValue "E40" "1.82%"
Set-TextValue "E42" "0.17%"
Set-TextValue "E43" "0.11%"
Set-TextValue "D44" "0.01388"
Set-TextValue "E44" "-15.23%"
Set-TextValue "E45" "-3.10%"
Set-TextValue "E46" "0.01%"
Set-TextValue "D48" "0.2503"
Set-TextValue "E48" "84.97%"
Set-TextValue "E49" "0.01%"
Set-TextValue "E50" "0.01%"
